# Updated Employee Scorecards for September
# Fill in the July (C), August (D) and September (E) attendance figures
# for the "Overall Performance / Punctuality" table (rows 28-35), and
# move the active selection to E32 to match the author's final cursor
# position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter constant, used below to match the center-aligned style that
# Excel applies to these numeric cells (style index 2 in styles.xml).
$xlCenter = -4108

function Set-Attendance {
    param($row, $july, $august, $september)

    $cJuly = $ws.Cells.Item($row, 3)   # column C
    $cAug  = $ws.Cells.Item($row, 4)   # column D
    $cSep  = $ws.Cells.Item($row, 5)   # column E

    $cJuly.HorizontalAlignment = $xlCenter
    $cAug.HorizontalAlignment  = $xlCenter
    $cSep.HorizontalAlignment  = $xlCenter

    $cJuly.Value = $july
    $cAug.Value  = $august
    $cSep.Value  = $september
}

# Row 28: Total Working Days
Set-Attendance 28 21 22 19
# Row 29: Employee Working Days
Set-Attendance 29 21 20 19
# Row 30: Leaves
Set-Attendance 30 0 2 0
# Row 31: Days Worked over 8 hours
Set-Attendance 31 2 4 3
# Row 32: Days Worked under 8 hours
Set-Attendance 32 0 0 5
# Row 33: Weekends
Set-Attendance 33 0 2 0
# Row 34: Late Arrivals
Set-Attendance 34 0 0 0
# Row 35: Short Leaves
Set-Attendance 35 0 0 5

# Move the selection to match the saved workbook state.
[void]$ws.Range("E32").Select()
